# Reverse the order of the comma-separated "Recorded By" values in column G,
# except for rows whose value references admin@admin.com (those are left as-is).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G
    $val = $cell.Text

    if ($null -eq $val) { continue }
    if ($val -eq "") { continue }
    if ($val -notmatch ",") { continue }
    if ($val -match "admin@admin.com") { continue }

    $parts = $val -split ",\s*"
    $n = $parts.Count
    $reversed = @()
    for ($i = $n - 1; $i -ge 0; $i--) {
        $reversed += $parts[$i]
    }
    $newVal = [string]::Join(", ", $reversed)

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
